$wb = $excel.ActiveWorkbook
$ws = $null

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value2 = 447.89474; $ws.Range("I2").Value2 = 423.84616; $ws.Range("K2").Value2 = 423.84616; $ws.Range("M2").Value2 = -310.84616
$ws.Range("H48").Value2 = 12748.167; $ws.Range("J48").Value2 = 12748.167; $ws.Range("L48").Value2 = 38244.501; $ws.Range("N48").Value2 = -38828.501
$ws.Range("H56").Value2 = 12748.167; $ws.Range("J56").Value2 = 12748.167; $ws.Range("L56").Value2 = 38244.501; $ws.Range("N56").Value2 = -39312.501
$ws.Range("H64").Value2 = 5788.8945; $ws.Range("I64").Value2 = 5561.8125; $ws.Range("J64").Value2 = 7000; $ws.Range("K64").Value2 = 5561.8125; $ws.Range("L64").Value2 = 7000; $ws.Range("M64").Value2 = -5313.8125; $ws.Range("N64").Value2 = -7496
$ws.Range("H67").Value2 = 5788.8945; $ws.Range("I67").Value2 = 5561.8125; $ws.Range("J67").Value2 = 7000; $ws.Range("K67").Value2 = 5561.8125; $ws.Range("L67").Value2 = 7000; $ws.Range("M67").Value2 = -4703.8125; $ws.Range("N67").Value2 = -8716
$ws.Range("H111").Value2 = 863.06665; $ws.Range("I111").Value2 = 795.5; $ws.Range("K111").Value2 = 2386.5; $ws.Range("M111").Value2 = 680.5
$ws.Range("H125").Value2 = 4469.8887; $ws.Range("I125").Value2 = 3866; $ws.Range("J125").Value2 = 5224.75; $ws.Range("K125").Value2 = 34794; $ws.Range("L125").Value2 = 47022.75; $ws.Range("M125").Value2 = -32334; $ws.Range("N125").Value2 = -51942.75
$ws.Range("H137").Value2 = 521074; $ws.Range("I137").Value2 = 2235.625; $ws.Range("J137").Value2 = 728609.4; $ws.Range("K137").Value2 = 6706.875; $ws.Range("L137").Value2 = 2185828.2; $ws.Range("M137").Value2 = -4156.875; $ws.Range("N137").Value2 = -2190928.2
$ws.Range("H138").Value2 = 38715; $ws.Range("J138").Value2 = 3498.8333; $ws.Range("L138").Value2 = 10496.4999; $ws.Range("N138").Value2 = -20776.4999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value2 = 31500; $ws.Range("J24").Value2 = 31500; $ws.Range("L24").Value2 = 31500; $ws.Range("N24").Value2 = -32248
$ws.Range("H32").Value2 = 3852.1953; $ws.Range("I32").Value2 = 1652.8052; $ws.Range("K32").Value2 = 1652.8052; $ws.Range("M32").Value2 = -1365.8052
$ws.Range("H74").Value2 = 4177.826; $ws.Range("I74").Value2 = 2234.4; $ws.Range("K74").Value2 = 2234.4; $ws.Range("M74").Value2 = -1360.4
$ws.Range("H77").Value2 = 4177.826; $ws.Range("I77").Value2 = 2234.4; $ws.Range("K77").Value2 = 11172; $ws.Range("M77").Value2 = -6804
$ws.Range("H100").Value2 = 31500; $ws.Range("J100").Value2 = 31500; $ws.Range("L100").Value2 = 31500; $ws.Range("N100").Value2 = -33664
$ws.Range("H122").Value2 = 2069.25; $ws.Range("I122").Value2 = 2138.2; $ws.Range("J122").Value2 = 2020; $ws.Range("K122").Value2 = 6414.599999999999; $ws.Range("L122").Value2 = 6060; $ws.Range("M122").Value2 = -3964.599999999999; $ws.Range("N122").Value2 = -10960
$ws.Range("H132").Value2 = 3072.9546; $ws.Range("I132").Value2 = 3083.5789; $ws.Range("J132").Value2 = 3005.6667; $ws.Range("K132").Value2 = 9250.736699999999; $ws.Range("L132").Value2 = 9017.000100000001; $ws.Range("M132").Value2 = -6720.736699999999; $ws.Range("N132").Value2 = -14077.0001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value2 = 153.11111; $ws.Range("I7").Value2 = 55.666668; $ws.Range("J7").Value2 = 222.71428; $ws.Range("K7").Value2 = 55.666668; $ws.Range("L7").Value2 = 222.71428; $ws.Range("M7").Value2 = 57.333332; $ws.Range("N7").Value2 = -448.71428
$ws.Range("H15").Value2 = 2724.1667; $ws.Range("J15").Value2 = 3455.7144; $ws.Range("L15").Value2 = 3455.7144; $ws.Range("N15").Value2 = -3795.7144
$ws.Range("H31").Value2 = 2660.5; $ws.Range("I31").Value2 = 1661; $ws.Range("J31").Value2 = 4104.222; $ws.Range("K31").Value2 = 1661; $ws.Range("L31").Value2 = 4104.222; $ws.Range("M31").Value2 = -1366; $ws.Range("N31").Value2 = -4694.222
$ws.Range("H34").Value2 = 2660.5; $ws.Range("I34").Value2 = 1661; $ws.Range("J34").Value2 = 4104.222; $ws.Range("K34").Value2 = 1661; $ws.Range("L34").Value2 = 4104.222; $ws.Range("M34").Value2 = -1459; $ws.Range("N34").Value2 = -4508.222
$ws.Range("H58").Value2 = 4451.923; $ws.Range("I58").Value2 = 4901.615; $ws.Range("J58").Value2 = 4002.2307; $ws.Range("K58").Value2 = 4901.615; $ws.Range("L58").Value2 = 4002.2307; $ws.Range("M58").Value2 = -4698.615; $ws.Range("N58").Value2 = -4408.2307
$ws.Range("H99").Value2 = 15875516; $ws.Range("I99").Value2 = 27779278; $ws.Range("J99").Value2 = 3833.3333; $ws.Range("K99").Value2 = 27779278; $ws.Range("L99").Value2 = 3833.3333; $ws.Range("M99").Value2 = -27777780; $ws.Range("N99").Value2 = -6829.3333
$ws.Range("H126").Value2 = 15875516; $ws.Range("I126").Value2 = 27779278; $ws.Range("J126").Value2 = 3833.3333; $ws.Range("K126").Value2 = 83337834; $ws.Range("L126").Value2 = 11499.9999; $ws.Range("M126").Value2 = -83335364; $ws.Range("N126").Value2 = -16439.9999
$ws.Range("H132").Value2 = 1174415.5; $ws.Range("I132").Value2 = 1202045.5; $ws.Range("K132").Value2 = 3606136.5; $ws.Range("M132").Value2 = -3603606.5
$ws.Range("H134").Value2 = 4831755; $ws.Range("I134").Value2 = 7146469; $ws.Range("J134").Value2 = 202326.8; $ws.Range("K134").Value2 = 21439407; $ws.Range("L134").Value2 = 606980.3999999999; $ws.Range("M134").Value2 = -21436872; $ws.Range("N134").Value2 = -612050.3999999999
$ws.Range("H136").Value2 = 4451.923; $ws.Range("I136").Value2 = 4901.615; $ws.Range("J136").Value2 = 4002.2307; $ws.Range("K136").Value2 = 14704.845; $ws.Range("L136").Value2 = 12006.6921; $ws.Range("M136").Value2 = -12154.845; $ws.Range("N136").Value2 = -17106.6921
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value2 = 10.107142; $ws.Range("I38").Value2 = 9.625; $ws.Range("J38").Value2 = 13; $ws.Range("K38").Value2 = 28.875; $ws.Range("L38").Value2 = 39; $ws.Range("M38").Value2 = 318.125; $ws.Range("N38").Value2 = -733
$ws.Range("H68").Value2 = 1270.5714; $ws.Range("J68").Value2 = 1211.25; $ws.Range("L68").Value2 = 3633.75; $ws.Range("N68").Value2 = -5255.75
$ws.Range("H71").Value2 = 1270.5714; $ws.Range("J71").Value2 = 1211.25; $ws.Range("L71").Value2 = 10901.25; $ws.Range("N71").Value2 = -19013.25
$ws.Range("H86").Value2 = 724.5; $ws.Range("I86").Value2 = 499; $ws.Range("J86").Value2 = 950; $ws.Range("K86").Value2 = 1497; $ws.Range("L86").Value2 = 2850; $ws.Range("M86").Value2 = -311; $ws.Range("N86").Value2 = -5222
$ws.Range("H89").Value2 = 724.5; $ws.Range("I89").Value2 = 499; $ws.Range("J89").Value2 = 950; $ws.Range("K89").Value2 = 4491; $ws.Range("L89").Value2 = 8550; $ws.Range("M89").Value2 = 1437; $ws.Range("N89").Value2 = -20406
$ws.Range("H132").Value2 = 6727.5557; $ws.Range("J132").Value2 = 6727.5557; $ws.Range("L132").Value2 = 60548.0013; $ws.Range("N132").Value2 = -65608.0013
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value2 = 262.94116; $ws.Range("I2").Value2 = 202.41667; $ws.Range("K2").Value2 = 202.41667; $ws.Range("M2").Value2 = -89.41667000000001
$ws.Range("H33").Value2 = 11742952; $ws.Range("J33").Value2 = 14550428; $ws.Range("L33").Value2 = 14550428; $ws.Range("N33").Value2 = -14550932
$ws.Range("H102").Value2 = 1540.4348; $ws.Range("I102").Value2 = 1412.8823; $ws.Range("K102").Value2 = 1412.8823; $ws.Range("M102").Value2 = 209.1177
$ws.Range("H122").Value2 = 12851.4; $ws.Range("I122").Value2 = 11814.25; $ws.Range("J122").Value2 = 17000; $ws.Range("K122").Value2 = 35442.75; $ws.Range("L122").Value2 = 51000; $ws.Range("M122").Value2 = -32992.75; $ws.Range("N122").Value2 = -55900
$ws.Range("H132").Value2 = 4533.7896; $ws.Range("I132").Value2 = 3301.0303; $ws.Range("J132").Value2 = 12670; $ws.Range("K132").Value2 = 9903.090899999999; $ws.Range("L132").Value2 = 38010; $ws.Range("M132").Value2 = -7373.090899999999; $ws.Range("N132").Value2 = -43070
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value2 = 2369.12; $ws.Range("I16").Value2 = 2048.8; $ws.Range("J16").Value2 = 3650.4; $ws.Range("K16").Value2 = 2048.8; $ws.Range("L16").Value2 = 3650.4; $ws.Range("M16").Value2 = -1878.8; $ws.Range("N16").Value2 = -3990.4
$ws.Range("H68").Value2 = 3033.3333; $ws.Range("J68").Value2 = 2807; $ws.Range("L68").Value2 = 2807; $ws.Range("N68").Value2 = -4305
$ws.Range("H71").Value2 = 3033.3333; $ws.Range("J71").Value2 = 2807; $ws.Range("L71").Value2 = 14035; $ws.Range("N71").Value2 = -21523
$ws.Range("H104").Value2 = 15000; $ws.Range("J104").Value2 = 15000; $ws.Range("L104").Value2 = 15000; $ws.Range("N104").Value2 = -21988
$ws.Range("H129").Value2 = 63193; $ws.Range("J129").Value2 = 49996; $ws.Range("L129").Value2 = 49996; $ws.Range("N129").Value2 = -59996
$ws.Range("H132").Value2 = 3915.4211; $ws.Range("I132").Value2 = 3770.1765; $ws.Range("J132").Value2 = 5150; $ws.Range("K132").Value2 = 11310.5295; $ws.Range("L132").Value2 = 15450; $ws.Range("M132").Value2 = -8780.529500000001; $ws.Range("N132").Value2 = -20510
$ws.Range("H141").Value2 = 97664.75; $ws.Range("J141").Value2 = 97664.75; $ws.Range("L141").Value2 = 97664.75; $ws.Range("N141").Value2 = -108024.75
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value2 = 4100.731; $ws.Range("I122").Value2 = 4256.9546; $ws.Range("K122").Value2 = 12770.8638; $ws.Range("M122").Value2 = -10320.8638
$ws.Range("H136").Value2 = 2549.7354; $ws.Range("I136").Value2 = 2475.4849; $ws.Range("K136").Value2 = 7426.4547; $ws.Range("M136").Value2 = -4876.4547